# Automatische test-sync: 2025-06-29 14:35:50
# Appends a new mail-log entry (row 14) to the "Logs" sheet, extends the
# conditional formatting ranges to cover the new row, and bumps the
# "Bestelling / Levering" counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append new row 14 to the Logs sheet ---------------------------------
$row = 14

$logs.Cells.Item($row, 1).Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$logs.Cells.Item($row, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($row, 5).Value = "Geachte klant,`r`n`r`nBedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.`r`n`r`nIk zie uw reactie graag tegemoet.`r`n`r`nMet vriendelijke groet,`r`n`r`n[Naam]`r`nE-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-29 14:35:08"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Ja"
$logs.Cells.Item($row, 9).Value = "Nee"

# Re-fit the row height after inserting multi-line text so no stray
# explicit row height is persisted (matches the other data rows).
$logs.Rows.Item($row).AutoFit()

# --- Extend conditional formatting ranges from row 13 to row 14 ----------
$logs.Range("D2:D13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D14"))
$logs.Range("G2:G13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G14"))
$logs.Range("H2:H13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H14"))
$logs.Range("I2:I13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I14"))

# --- Update Dashboard summary count for "Bestelling / Levering" ----------
$dash.Range("B3").Value = 4
